$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, copying the style used by the other headers (e.g. G1)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill H2:H10 with 0 (numeric)
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
